$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows (row 2..5) with new "preguntador N / pregunta N / respuesta N" values.
# Order follows the author's original edit order (deduced from the resulting
# shared-string table order in the target file).
$ws.Range("A2").Value = "preguntador 1"
$ws.Range("B2").Value = "pregunta 1"
$ws.Range("C2").Value = "respuesta 1"

$ws.Range("C3").Value = "respuesta 2"
$ws.Range("C4").Value = "respuesta 3"
$ws.Range("C5").Value = "respuesta 4"

$ws.Range("A3").Value = "preguntador 2"
$ws.Range("A4").Value = "preguntador 3"
$ws.Range("A5").Value = "preguntador 4"

$ws.Range("B3").Value = "pregunta 2"
$ws.Range("B4").Value = "pregunta 3"
$ws.Range("B5").Value = "pregunta 4"

# Resize columns B and C (author narrowed them after editing the text)
$ws.Columns("B").ColumnWidth = 13
$ws.Columns("C").ColumnWidth = 12.42578125

# --- style bookkeeping -----------------------------------------------
# C3 previously carried a one-off underlined style; drop it so it matches
# the plain text style shared by the rest of column C.
$g5 = $ws.Range("G5")
$g5.NumberFormat = "@"
$g5.Font.Underline = $true

$c3 = $ws.Range("C3")
$c3.Font.Underline = $false

$g5.NumberFormatLocal = "General"

# Click on G5 (empty cell outside the table) - this is the final selection
$g5.Select()
